# Generate Report for Handback
#
# Marks the two tracked files ("7a0dda7e-...md" and "dd96ed0d-...md") as
# handed back (in sync with en-US) for both locales (zh-cn, de-de):
#   - Status column (C) text updated
#   - Latest Target File (F) / Latest Handback File (G) populated with
#     hyperlinked file names
#   - Latest Handback DateTime (H) stamped with the handback timestamp

$wb = $excel.ActiveWorkbook

$mdHash        = "bfb1c9ad15365171df98d65adec4de54589aae2a"
$statusText    = "Handed back: in sync with en-US"

$file1Base     = "7a0dda7e-6cfd-4672-9cd0-e43464cc48aa"
$file1Hash     = "fbf7016bae921db29f78584cb580c6ddfafb0374"
$file2Base     = "dd96ed0d-1caa-4160-84fb-cd10959e7297"
$file2Hash     = "9a64b914d41de930e717736ca7c02df57a1b744e"

# Per-locale configuration: sheet name, handoff-commit hash (for the
# "ht" target link - matches what's already on column D), handback
# timestamp written to column H.
$locales = @(
    @{ Sheet = "zh-cn"; Code = "zh-cn"; HoHash = "5c7648591d59bae75aa4472f4e0a0f9022eff20b"; Stamp = "2016-03-19 17:12:32"; FlyDir = "oltest-zhcn-fly" },
    @{ Sheet = "de-de"; Code = "de-de"; HoHash = "9cdd50f5a67af2db7f0bfc8d09b47a6fa2062294"; Stamp = "2016-03-19 17:12:45"; FlyDir = "oltest-dede-fly" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # --- Status text (row 2 and row 3) ---
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # --- Row 2 (file1 = 7a0dda7e...) ---
    $f1Md   = "$file1Base.md"
    $f1Xlf  = "$file1Base.$file1Hash.$($loc.Code).xlf"

    $ws.Range("F2").Value = $f1Md
    $ws.Hyperlinks.Add(
        $ws.Range("F2"),
        "https://github.com/OpenLocalizationTest/oltest/blob/$mdHash/e2e/$f1Md",
        "",
        "",
        $f1Md
    )

    $ws.Range("G2").Value = $f1Xlf
    $ws.Hyperlinks.Add(
        $ws.Range("G2"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($loc.HoHash)/ol-handback/OpenLocalizationTest/oltest/yuwzho/hb/$f1Xlf",
        "",
        "",
        $f1Xlf
    )

    $ws.Range("H2").Value = $loc.Stamp

    # --- Row 3 (file2 = dd96ed0d...) ---
    $f2Md   = "$file2Base.md"
    $f2Xlf  = "$file2Base.$file2Hash.$($loc.Code).xlf"

    $ws.Range("F3").Value = $f2Md
    $ws.Hyperlinks.Add(
        $ws.Range("F3"),
        "https://github.com/OpenLocalizationTest/oltest/blob/$mdHash/e2e/$f2Md",
        "",
        "",
        $f2Md
    )

    $ws.Range("G3").Value = $f2Xlf
    $ws.Hyperlinks.Add(
        $ws.Range("G3"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($loc.HoHash)/ol-handback/OpenLocalizationTest/oltest/yuwzho/hb/$f2Xlf",
        "",
        "",
        $f2Xlf
    )

    $ws.Range("H3").Value = $loc.Stamp
}
